$p = $ppt.ActivePresentation

# --- Update the "Action.DisplayMessageForm/Action.DisplayAppointmentForm" box on
# --- slide 5 ("Outlook-specific Adaptive Card properties and features" / ActionSet)
# --- to describe the new Action.Transaction action instead.
$s5 = $p.Slides.Item(5)
$sh = $s5.Shapes.Item(4)
$tr = $sh.TextFrame.TextRange
$tr.Paragraphs(1).Runs(1).Text = "Action.Transaction"
$tr.Paragraphs(2).Delete()
$tr.Paragraphs(2).Runs(1).Text = "Triggers the payments in Outlook experience using Microsoft Pay."
$sh.Width = 298.1118927637795
$sh.Height = 89.53590781181103

# --- Remove the "Actions" summary table slide (Action/Description table with
# --- OpenUri / HttpPOST / ActionCard rows); its content was folded elsewhere.
$p.Slides.Item(6).Delete()
